$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") for rows 2-27 was bumped by one day:
# serial 45316 -> 45317 (2024-01-25 -> 2024-01-26).
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45316) {
        $cell.Value2 = 45317
    }
}
